$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Total
$ws.Range("A9").Value = "Total"
$ws.Range("B9").Formula = "=SUM(B2:B6)"
$ws.Range("C9").Formula = "=SUM(C2:C6)"
$ws.Range("D9").Formula = "=SUM(D2:D6)"
$ws.Range("E9").Formula = "=C9/B9"
$ws.Range("F9").Formula = "=D9/B9"

# Row 10: Neither over nor under 95CI
$ws.Range("A10").Value = "Neither over nor under 95CI"
$ws.Range("B10").Formula = "=B9-(C9+D9)"

# Row 11: pct not over/under
$ws.Range("B11").Formula = "=B10/B9"

# Formatting: rows 9-11 use the red "Tipping Point" font (same as row header E1/F1, fontId 14)
$ws.Range("A9:D9").Font.Color = 255
$ws.Range("A10:D10").Font.Color = 255
$ws.Range("A11").Font.Color = 255
$ws.Range("C11:D11").Font.Color = 255
$ws.Range("B11").Font.Color = 255

# Percent number formats for ratio columns
$ws.Range("E9:F9").NumberFormat = "0%"
$ws.Range("B11").NumberFormat = "0%"

$ws.Range("D10").Select()
